$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# List of cell updates: (cell address, new text value)
# Values are written as plain text to match the source data (which stores
# prices/percentages as literal strings, not numbers), by temporarily
# forcing a Text number format, then clearing formatting again so the
# cell's style stays at the workbook default.
$updates = @(
    @{Cell="D2"; Value="65.845.51"},
    @{Cell="E2"; Value="  +0.56%  "},
    @{Cell="D3"; Value="2.680.68"},
    @{Cell="E3"; Value="  +0.74%  "},
    @{Cell="E4"; Value="  +0.03%  "},
    @{Cell="D5"; Value="603.19"},
    @{Cell="E5"; Value="  -0.37%  "},
    @{Cell="D6"; Value="156.75"},
    @{Cell="E6"; Value="  -0.62%  "},
    @{Cell="E7"; Value="  +0.02%  "},
    @{Cell="D8"; Value="0.589"},
    @{Cell="E8"; Value="  -0.09%  "},
    @{Cell="E9"; Value="  -0.05%  "},
    @{Cell="D10"; Value="5.94"},
    @{Cell="E10"; Value="  +2.00%  "},
    @{Cell="E11"; Value="  -3.22%  "},
    @{Cell="E12"; Value="  +0.22%  "},
    @{Cell="D13"; Value="29.50"},
    @{Cell="E13"; Value="  -0.38%  "},
    @{Cell="D14"; Value="0.0000200"},
    @{Cell="E14"; Value="  +5.88%  "},
    @{Cell="D15"; Value="3.162.49"},
    @{Cell="E15"; Value="  +0.74%  "},
    @{Cell="D16"; Value="65.629.41"},
    @{Cell="E16"; Value="  +0.59%  "},
    @{Cell="D17"; Value="2.676.20"},
    @{Cell="E17"; Value="  +0.87%  "},
    @{Cell="D18"; Value="12.64"},
    @{Cell="E18"; Value="  -1.40%  "},
    @{Cell="E19"; Value="  -1.79%  "},
    @{Cell="D20"; Value="7.60"},
    @{Cell="E20"; Value="  +3.02%  "},
    @{Cell="D21"; Value="352.91"},
    @{Cell="E21"; Value="  -2.00%  "},
    @{Cell="E22"; Value="  -0.15%  "},
    @{Cell="D23"; Value="70.37"},
    @{Cell="E23"; Value="  +1.79%  "},
    @{Cell="E24"; Value="  +7.22%  "},
    @{Cell="D25"; Value="9.87"},
    @{Cell="E25"; Value="  +3.41%  "},
    @{Cell="D26"; Value="1.62"},
    @{Cell="E26"; Value="  -5.60%  "},
    @{Cell="E27"; Value="  -1.65%  "},
    @{Cell="D28"; Value="0.169"},
    @{Cell="E28"; Value="  +1.93%  "},
    @{Cell="E29"; Value="  -1.25%  "},
    @{Cell="E30"; Value="  +0.33%  "},
    @{Cell="B31"; Value="Bittensor"},
    @{Cell="C31"; Value="https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"},
    @{Cell="D31"; Value="536.28"},
    @{Cell="E31"; Value="  -1.34%  "},
    @{Cell="B32"; Value="PancakeSwap"},
    @{Cell="C32"; Value="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"},
    @{Cell="D32"; Value="2.15"},
    @{Cell="E32"; Value="  -3.79%  "},
    @{Cell="E33"; Value="  -4.25%  "},
    @{Cell="E34"; Value="  +2.07%  "},
    @{Cell="D35"; Value="5.41"},
    @{Cell="E35"; Value="  -4.58%  "},
    @{Cell="D36"; Value="0.427"},
    @{Cell="E36"; Value="  -1.41%  "},
    @{Cell="D37"; Value="20.49"},
    @{Cell="E37"; Value="  -1.06%  "},
    @{Cell="D38"; Value="160.69"},
    @{Cell="E38"; Value="  -1.76%  "},
    @{Cell="E39"; Value="  -1.87%  "},
    @{Cell="E40"; Value="  -0.01%  "},
    @{Cell="D41"; Value="1.00"},
    @{Cell="D42"; Value="42.29"},
    @{Cell="E42"; Value="  -0.64%  "},
    @{Cell="D43"; Value="166.33"},
    @{Cell="E43"; Value="  -0.42%  "},
    @{Cell="E44"; Value="  -2.24%  "},
    @{Cell="D45"; Value="0.0621"},
    @{Cell="E45"; Value="  +0.13%  "},
    @{Cell="D46"; Value="23.15"},
    @{Cell="E46"; Value="  -0.05%  "},
    @{Cell="E47"; Value="  -4.44%  "},
    @{Cell="E48"; Value="  -0.31%  "},
    @{Cell="D49"; Value="0.651"},
    @{Cell="E49"; Value="  -1.37%  "},
    @{Cell="E50"; Value="  +2.58%  "},
    @{Cell="D51"; Value="0.0988"},
    @{Cell="E51"; Value="  +0.16%  "}
)

foreach ($item in $updates) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
    $rng.ClearFormats()
}
